# Apply "feat: integrations with sheel" edit:
#  - entidade: update cnpj (C2) and inscricao_estadual (E2), select C6, make it the active/visible tab
#  - ftp: append "2" to the ftp path cells in row 2, move selection to A19
#  - tipo_recebimento: loses tab-selected state (handled automatically by activating entidade)

$wb = $excel.ActiveWorkbook

# --- Update "ftp" sheet values (importacao/exportacao/erro paths get a trailing "2") ---
$wsFtp = $wb.Worksheets.Item("ftp")
$wsFtp.Range("A2").Value = "/home/ftpsynapcomp/Embu/Vtex/importacao232"
$wsFtp.Range("B2").Value = "/home/ftpsynapcomp/Embu/Vtex/bkp_importacao232"
$wsFtp.Range("C2").Value = "/home/ftpsynapcomp/Embu/Vtex/exportacao232"
$wsFtp.Range("D2").Value = "/home/ftpsynapcomp/Embu/Vtex/bkp_exportacao232"
$wsFtp.Range("E2").Value = "/home/ftpsynapcomp/Embu/Vtex/erro232"

# --- Update "entidade" sheet values (new cnpj / inscricao_estadual) ---
$wsEntidade = $wb.Worksheets.Item("entidade")
$wsEntidade.Range("C2").Value = "07.046.881/1007-14"
$wsEntidade.Range("E2").Value = 15368

# Move selection within "ftp" before leaving it, so it is saved as the sheet's selection
$wsFtp.Activate()
$wsFtp.Range("A19").Select()

# Activate "entidade" and select C6 - this becomes the selected/visible tab
$wsEntidade.Activate()
$wsEntidade.Range("C6").Select()
